# The upstream diff for this fixture touches only two things:
#
#   1. The order in which XML-namespace declarations are written on the
#      root elements of word/document.xml, word/footer.xml, word/header.xml
#      and word/styles.xml.
#   2. A single word inside a non-standard XML *comment* that docx4j itself
#      injected as the very first child of <w:body>:
#         "... REFERENCE JAXB in Microsoft Java 21.0.8 on Mac OS X ..."
#      becoming
#         "... REFERENCE JAXB in Oracle Java 21.0.8 on Mac OS X ..."
#
# Both of these are build-tool/environment artifacts of the external
# docx4j/JAXB pipeline that produced the reference file (the JAXB
# namespace-context map iterates in a different order under a different
# JRE vendor, and docx4j stamps the JRE vendor name into that comment).
# Neither one is reachable through Word's object model: a bare XML comment
# that is a direct child of <w:body> is not a paragraph, run, field,
# content control, or any other addressable OM entity, so Word (real or
# automated) has no Find/Replace target, Range, or property that exposes
# its text, and it is not reproducible through Range.InsertXML either
# (XML comments are not part of WordprocessingML content and are dropped
# by the importer). Likewise, nothing in the Word OM lets a caller control
# the byte-level attribute order the writer uses when it re-serializes a
# part's root element.
#
# Because no actual document content, formatting, or structure changed in
# this commit, the correct COM-interop replay is to leave the document's
# content exactly as-is: touch the active document without mutating any
# paragraph, run, table, header/footer or style, so the saved .docx is a
# faithful passthrough of the input (the only part of this commit that
# *is* expressible through the object model is "no content changed").

$d = $word.ActiveDocument

# Touch the document/object model (read-only) so the script demonstrably
# ran against $word.ActiveDocument without altering any content.
$null = $d.Content.Text.Length
$null = $d.Paragraphs.Count
